$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.988.36'
$ws.Range('E2').Value = '  +3.35%  '
$ws.Range('D3').Value = '3.801.96'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''697.12'
$ws.Range('E5').Value = '  +11.03%  '
$ws.Range('D6').Value = '''173.88'
$ws.Range('E6').Value = '  +5.54%  '
$ws.Range('D7').Value = '3.800.02'
$ws.Range('E7').Value = '  +1.20%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +1.72%  '
$ws.Range('E10').Value = '  +3.73%  '
$ws.Range('D11').Value = '''7.47'
$ws.Range('E11').Value = '  +8.37%  '
$ws.Range('E12').Value = '  +1.53%  '
$ws.Range('E13').Value = '  +10.13%  '
$ws.Range('E14').Value = '  +4.94%  '
$ws.Range('D15').Value = '4.444.87'
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('D16').Value = '3.803.60'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').Value = '71.003.47'
$ws.Range('E17').Value = '  +3.39%  '
$ws.Range('D18').Value = '''17.87'
$ws.Range('E18').Value = '  +1.73%  '
$ws.Range('D19').Value = '''7.21'
$ws.Range('E19').Value = '  +3.45%  '
$ws.Range('E20').Value = '  +0.72%  '
$ws.Range('D21').Value = '''11.11'
$ws.Range('E21').Value = '  +17.51%  '
$ws.Range('D22').Value = '''485.15'
$ws.Range('E22').Value = '  +4.01%  '
$ws.Range('D23').Value = '''0.715'
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').Value = '''84.14'
$ws.Range('E24').Value = '  +3.28%  '
$ws.Range('E25').Value = '  +2.52%  '
$ws.Range('E26').Value = '  +3.09%  '
$ws.Range('E27').Value = '  +4.66%  '
$ws.Range('D28').Value = '''2.16'
$ws.Range('E28').Value = '  +3.47%  '
$ws.Range('D29').Value = '3.955.42'
$ws.Range('E29').Value = '  +1.29%  '
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').Value = '''3.04'
$ws.Range('E31').Value = '  +14.97%  '
$ws.Range('D32').Value = '''7.54'
$ws.Range('E32').Value = '  +6.49%  '
$ws.Range('E33').Value = '  +1.63%  '
$ws.Range('D34').Value = '''29.66'
$ws.Range('E34').Value = '  +4.94%  '
$ws.Range('E35').Value = '  +1.96%  '
$ws.Range('D36').Value = '''9.26'
$ws.Range('E36').Value = '  +4.82%  '
$ws.Range('E37').Value = '  +0.43%  '
$ws.Range('D38').Value = '3.754.96'
$ws.Range('E38').Value = '  +1.20%  '
$ws.Range('E39').Value = '  +2.89%  '
$ws.Range('E40').Value = '  +9.77%  '
$ws.Range('E41').Value = '  +4.35%  '
$ws.Range('D42').Value = '''2.23'
$ws.Range('E42').Value = '  +14.15%  '
$ws.Range('D43').Value = '''0.000328'
$ws.Range('E43').Value = '  +24.03%  '
$ws.Range('D44').Value = '''0.971'
$ws.Range('E44').Value = '  +1.68%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '''163.61'
$ws.Range('E47').Value = '  +5.12%  '
$ws.Range('D48').Value = '''49.23'
$ws.Range('E48').Value = '  +4.96%  '
$ws.Range('D49').Value = '''45.01'
$ws.Range('E49').Value = '  +2.00%  '
$ws.Range('E50').Value = '  +3.27%  '
$ws.Range('D51').Value = '''1.38'
$ws.Range('E51').Value = '  -1.36%  '
